$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume change %) refreshed from source feed.
$updates = @(
    @{ Row = 2; D = "28.780.62"; E = "  +2.57%  " }
    @{ Row = 3; D = "1.873.79"; E = "  +2.18%  " }
    @{ Row = 4; D = "1.005"; E = "  +0.37%  " }
    @{ Row = 5; D = "324.54"; E = "  +0.13%  " }
    @{ Row = 6; D = "1.004"; E = "  +0.32%  " }
    @{ Row = 7; D = "0.4619"; E = "  -0.45%  " }
    @{ Row = 8; D = $null; E = "  -0.25%  " }
    @{ Row = 9; D = "0.07865"; E = "  +0.12%  " }
    @{ Row = 10; D = "0.9881"; E = "  +2.69%  " }
    @{ Row = 11; D = "21.85"; E = "  -0.35%  " }
    @{ Row = 12; D = "1.897.29"; E = "  +2.19%  " }
    @{ Row = 13; D = "6.994"; E = "  +1.12%  " }
    @{ Row = 14; D = "5.705"; E = "  +0.30%  " }
    @{ Row = 15; D = "0.06977"; E = "  +2.04%  " }
    @{ Row = 16; D = "88.38"; E = "  +0.03%  " }
    @{ Row = 17; D = "1.005"; E = "  +0.37%  " }
    @{ Row = 18; D = "0.00001004"; E = "  +0.94%  " }
    @{ Row = 19; D = "16.80"; E = "  +0.54%  " }
    @{ Row = 20; D = $null; E = "  +0.18%  " }
    @{ Row = 21; D = "28.803.80"; E = "  +2.57%  " }
    @{ Row = 22; D = "5.283"; E = "  -0.55%  " }
    @{ Row = 23; D = "11.07"; E = "  +0.58%  " }
    @{ Row = 24; D = $null; E = "  +0.26%  " }
    @{ Row = 25; D = "2.129.93"; E = "  +2.98%  " }
    @{ Row = 26; D = "153.15"; E = "  -1.14%  " }
    @{ Row = 27; D = $null; E = "  +0.42%  " }
    @{ Row = 28; D = "5.851"; E = "  +3.06%  " }
    @{ Row = 29; D = "1.993"; E = "  +1.61%  " }
    @{ Row = 30; D = "118.94"; E = "  +0.57%  " }
    @{ Row = 31; D = "0.09334"; E = "  +1.08%  " }
    @{ Row = 32; D = "0.9198"; E = "  -1.75%  " }
    @{ Row = 33; D = "5.312"; E = "  +1.01%  " }
    @{ Row = 34; D = "1.338"; E = "  +1.24%  " }
    @{ Row = 35; D = "3.323"; E = "  +0.47%  " }
    @{ Row = 36; D = "0.05793"; E = "  -1.36%  " }
    @{ Row = 37; D = "1.151"; E = "  +0.44%  " }
    @{ Row = 38; D = "0.02072"; E = "  -2.61%  " }
    @{ Row = 39; D = "7.664"; E = "  -1.31%  " }
    @{ Row = 40; D = "0.5643"; E = "  +0.80%  " }
    @{ Row = 41; D = $null; E = "  +1.11%  " }
    @{ Row = 42; D = "9.808"; E = "  -0.88%  " }
    @{ Row = 43; D = "0.07212"; E = "  -0.89%  " }
    @{ Row = 44; D = "11.76"; E = "  +0.38%  " }
    @{ Row = 45; D = "0.5299"; E = "  +0.43%  " }
    @{ Row = 46; D = "2.127"; E = "  +0.65%  " }
    @{ Row = 47; D = "1.123"; E = "  -2.74%  " }
    @{ Row = 48; D = "1.838"; E = "  +0.66%  " }
    @{ Row = 49; D = "113.20"; E = "  +0.88%  " }
    @{ Row = 50; D = "2.419"; E = "  +3.83%  " }
    @{ Row = 51; D = $null; E = "  +0.27%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
